$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "License Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Delete()

# ---------------------------------------------------------------------------
# 2. Rebuild the big license paragraph (now paragraph 4) with the new text.
#    Clear its existing runs (leave the paragraph mark + the pre-existing
#    leading empty run in place) and insert the replacement content.
# ---------------------------------------------------------------------------
$bigPara = $d.Paragraphs.Item(4)
$paraStart = $bigPara.Range.Start
$paraEnd = $bigPara.Range.End

# Clear everything except the trailing paragraph mark.
$clearRng = $d.Range($paraStart, $paraEnd - 1)
$clearRng.Text = ""

$boldText = "unfoldingWord® Translation Questions"
$restText = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. unfoldingWord® Translation Questions has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from unfoldingWord® Translation Questions © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

$insPoint = $d.Range($paraStart, $paraStart)
$insPoint.InsertAfter($boldText + $restText)

$boldRng = $d.Range($paraStart, $paraStart + $boldText.Length)
$boldRng.Bold = 1

# ---------------------------------------------------------------------------
# 3. Remove the "This PDF version is provided under the same license."
#    paragraph entirely (it merges away, leaving the big paragraph's own
#    trailing empty run as the final run of the paragraph).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(5).Range.Delete()
